# Updated cryptos list on Wed Jun  5 13:25:45 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # The "Price" column stores plain-looking numbers as text (inlineStr) in
    # the source data. Forcing the number format to Text before the write
    # (and resetting the style back to Normal afterwards) keeps the cell a
    # string instead of letting Excel auto-coerce it to a number, while not
    # leaving any stray number-format style behind.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "71.065.13"
$ws.Range("E2").Value = "  +2.58%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.813.58"
$ws.Range("E3").Value = "  +1.00%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
Set-TextValue "D5" "699.65"
$ws.Range("E5").Value = "  +7.77%  "

# Row 6 - Solana
Set-TextValue "D6" "173.73"
$ws.Range("E6").Value = "  +4.66%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.812.03"
$ws.Range("E7").Value = "  +0.96%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - XRP
Set-TextValue "D9" "0.527"
$ws.Range("E9").Value = "  +0.44%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.163"
$ws.Range("E10").Value = "  +2.38%  "

# Row 11 - Toncoin
Set-TextValue "D11" "7.22"
$ws.Range("E11").Value = "  +4.55%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.65%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000259"
$ws.Range("E13").Value = "  +8.31%  "

# Row 14 - Avalanche
Set-TextValue "D14" "36.30"
$ws.Range("E14").Value = "  +3.78%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.453.72"
$ws.Range("E15").Value = "  +1.00%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.802.38"
$ws.Range("E16").Value = "  +1.04%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "71.006.21"
$ws.Range("E17").Value = "  +2.63%  "

# Row 18 - Chainlink
Set-TextValue "D18" "17.82"
$ws.Range("E18").Value = "  +0.12%  "

# Row 19
$ws.Range("E19").Value = "  +2.49%  "

# Row 20
$ws.Range("E20").Value = "  +0.24%  "

# Row 21
Set-TextValue "D21" "11.11"
$ws.Range("E21").Value = "  +15.91%  "

# Row 22
Set-TextValue "D22" "479.73"
$ws.Range("E22").Value = "  +2.50%  "

# Row 23
$ws.Range("E23").Value = "  +0.53%  "

# Row 24
Set-TextValue "D24" "83.85"
$ws.Range("E24").Value = "  +2.40%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "12.32"
$ws.Range("E26").Value = "  -0.18%  "

# Row 27 - was RenderToken, now Fetch.AI (rows 27/28 swapped)
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D27" "2.16"
$ws.Range("E27").Value = "  +2.59%  "

# Row 28 - was Fetch.AI, now RenderToken
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D28" "10.45"
$ws.Range("E28").Value = "  +0.61%  "

# Row 29 - WrappedeETH
Set-TextValue "D29" "3.963.61"
$ws.Range("E29").Value = "  +1.03%  "

# Row 30 - Dai
Set-TextValue "D30" "0.999"
$ws.Range("E30").Value = "  -0.12%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "3.12"
$ws.Range("E31").Value = "  +15.32%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +1.55%  "

# Row 33 - NEARProtocol
Set-TextValue "D33" "7.53"
$ws.Range("E33").Value = "  +5.12%  "

# Row 34 - Kaspa
Set-TextValue "D34" "0.186"
$ws.Range("E34").Value = "  +7.98%  "

# Row 35 - EthereumClassic
Set-TextValue "D35" "29.53"
$ws.Range("E35").Value = "  +2.98%  "

# Row 36 - Aptos
Set-TextValue "D36" "9.26"
$ws.Range("E36").Value = "  +4.86%  "

# Row 37 - Binance-PegBSC-USD
Set-TextValue "D37" "0.999"
$ws.Range("E37").Value = "  -0.03%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  +2.71%  "

# Row 39 - dogwifhat
Set-TextValue "D39" "3.44"
$ws.Range("E39").Value = "  +5.76%  "

# Row 40 - Filecoin
Set-TextValue "D40" "6.01"
$ws.Range("E40").Value = "  +2.96%  "

# Row 41 - Stacks
Set-TextValue "D41" "2.25"
$ws.Range("E41").Value = "  +13.28%  "

# Row 42 - Mantle
Set-TextValue "D42" "0.978"
$ws.Range("E42").Value = "  +2.34%  "

# Row 43 - FLOKI
Set-TextValue "D43" "0.000327"
$ws.Range("E43").Value = "  +21.29%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  +0.02%  "

# Row 45 - USDe
$ws.Range("E45").Value = "  -0.02%  "

# Row 46 - Monero
Set-TextValue "D46" "163.15"
$ws.Range("E46").Value = "  +4.50%  "

# Row 47 - OKB
Set-TextValue "D47" "48.96"
$ws.Range("E47").Value = "  +3.50%  "

# Row 48 - Arweave
Set-TextValue "D48" "44.40"
$ws.Range("E48").Value = "  -2.04%  "

# Row 49 - TheGraph
$ws.Range("E49").Value = "  +1.39%  "

# Row 50 - ONDO
$ws.Range("E50").Value = "  -1.46%  "

# Row 51 - Bittensor
Set-TextValue "D51" "409.45"
$ws.Range("E51").Value = "  +6.62%  "
